# NE_chartering_2021_Aug_Borek_daily_totals.xlsx
#
# - costing of crew and weather days
# - add "fly_cost" column (after fly_time)
# - add "day cost per all people" column (before last activity)
# - add two weather-delay days (2021-08-09 and 2021-08-12)
# - add "grand total (MDKK)" and "grand total incl. quarantine (MDKK)" rows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Insert the two new "wx delay day" rows (while the sheet still has
#    its original 5 columns, so row numbers line up with the source data)
# ---------------------------------------------------------------------
# current: 2:08-07 3:08-08 4:08-10 5:08-11 6:08-13 7:08-14 8:total
$ws.Rows.Item(4).Insert()
# now:     2:08-07 3:08-08 4:(new) 5:08-10 6:08-11 7:08-13 8:08-14 9:total
$ws.Rows.Item(7).Insert()
# now:     2:08-07 3:08-08 4:(new) 5:08-10 6:08-11 7:(new) 8:08-13 9:08-14 10:total

# ---------------------------------------------------------------------
# 2) Insert the two new columns: "fly_cost" (C) and
#    "day cost per all people" (F, i.e. before "last activity")
# ---------------------------------------------------------------------
$ws.Columns.Item(3).Insert()
$ws.Columns.Item(6).Insert()

# ---------------------------------------------------------------------
# 3) Header row
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "date (YYYY-MM-DD)"
$ws.Range("B1").Value = "fly_time"
$ws.Range("C1").Value = "fly_cost"
$ws.Range("D1").Value = "fuel consumption litres"
$ws.Range("E1").Value = "fuel consumption cost kDKK"
$ws.Range("F1").Value = "day cost per all people"
$ws.Range("G1").Value = "last activity"

# carry over the bordered/bold header style onto the two new header cells
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4) Data rows (date column forced to text so it keeps its literal
#    "YYYY-MM-DD" form instead of becoming a date serial)
# ---------------------------------------------------------------------
$data = @(
    @(2,  "2021-08-07", 3.5, 46.6, 1391, 14.1, 7.5,   "transit from Iqaluit with no cargo"),
    @(3,  "2021-08-08", 4.3, 58.1, 2933, 29.7, 15,    "EGP new AWS install. Twin Otter overnights"),
    @(4,  "2021-08-09", 4,   53.6, 0,    0,    15,    "wx delay day"),
    @(5,  "2021-08-10", 1.4, 18.9, 846,  8.6,  15,    "return to EGP for overnight 2"),
    @(6,  "2021-08-11", 2.6, 35.4, 1742, 17.7, 15,    "return to coast. Upernavik"),
    @(7,  "2021-08-12", 4,   53.6, 0,    0,    15,    "wx delay day"),
    @(8,  "2021-08-13", 2.6, 35.4, 1731, 17.5, 15,    "refuel in JAV then end day in SFJ"),
    @(9,  "2021-08-14", 3.5, 46.6, 1391, 14.1, 15,    "return to Canada")
)

foreach ($row in $data) {
    $r = $row[0]
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $row[1]
    $dateCell.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
}

# ---------------------------------------------------------------------
# 5) "total" row
# ---------------------------------------------------------------------
$ws.Range("A10").Value = "total"
$ws.Range("B10").Value = 25.9
$ws.Range("C10").Value = 348.2
$ws.Range("D10").Value = 10034
$ws.Range("E10").Value = 101.7
$ws.Range("F10").Value = 112.5

# ---------------------------------------------------------------------
# 6) grand-total summary rows
# ---------------------------------------------------------------------
$ws.Range("A11").Value = "grand total (MDKK)"
$ws.Range("B11").Value = 0.5624000000000001

$ws.Range("A12").Value = "grand total incl. quarantine (MDKK)"
$ws.Range("B12").Value = 0.6374000000000001
